# Applies the cryptos-list refresh captured in the commit:
#   "Updated cryptos list on Sat Apr  6 03:13:12 UTC 2024 with GitHub Actions"
#
# All Coin/Link/Price/Volume cells are stored as plain text in the workbook
# (t="inlineStr"), so for any replacement price that Excel would otherwise
# auto-parse as a number (dropping meaningful trailing zeros, e.g. "1.00"
# becoming 1), the cells NumberFormat is pinned to "@" (Text) first -
# scoped to just that single cell so no other formatting is disturbed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.656.27"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.325.73"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.66"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.03"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "3.323.15"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.01"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "703.52"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "3.868.72"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.43"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "67.674.38"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "3.321.93"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.892"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.89"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.62"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.39"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.19"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.52"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  +5.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.50"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.96"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "3.687.15"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.68"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.13"
$ws.Range("E39").Value = "  +6.22%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "0.0₃0671"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.335"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +6.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.97"
$ws.Range("E51").Value = "  -0.10%  "
